# Remove the old "MGCycle" heading stub (and its surrounding blank
# paragraphs) that was left at the end of the document.
#
# The tail of the document currently looks like:
#   ...u=SOR(u,f)...                      <- last real content paragraph
#   <empty paragraph with eastAsiaTheme rFonts>
#   "MGCycle"   (Heading1)
#   <empty paragraph>
#   <empty paragraph>
#   <empty paragraph>
#   <empty paragraph>
#   <empty paragraph>
#
# All seven of the trailing paragraphs (the blank spacer, the MGCycle
# heading, and the five blank paragraphs after it) need to be deleted,
# so the document again ends right after the "u=SOR(u,f)" paragraph.

$d = $word.ActiveDocument

# Locate the "MGCycle" heading paragraph by its text so this doesn't
# depend on a hard-coded paragraph index.
$searchRange = $d.Content
$found = $searchRange.Find.Execute("MGCycle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $headingPara = $searchRange.Paragraphs.Item(1)
    $headingIndex = $headingPara.Index

    # The blank paragraph immediately before the heading is also part of
    # the block being removed.
    $firstParaToRemove = $d.Paragraphs.Item($headingIndex - 1)

    $deleteStart = $firstParaToRemove.Range.Start
    $deleteEnd = $d.Content.End

    $deleteRange = $d.Range($deleteStart, $deleteEnd)
    $deleteRange.Delete()
}
